# interactive graph: year and woc
# Update the id values in column A for rows 19-28 on the active sheet ("kosong").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 3339
$ws.Range("A20").Value = 3340
$ws.Range("A21").Value = 4582
$ws.Range("A22").Value = 4583
$ws.Range("A23").Value = 4584
$ws.Range("A24").Value = 6278
$ws.Range("A25").Value = 6279
$ws.Range("A26").Value = 6280
$ws.Range("A27").Value = 7634
$ws.Range("A28").Value = 7635
